$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections: cells that flip between numeric 1/value and text "NaN" ---
$ws.Range("BS13").Value = "NaN"
$ws.Range("AD25").Value = 1
$ws.Range("AD26").Value = 1
$ws.Range("AD27").Value = 1
$ws.Range("CM30").Value = 1
$ws.Range("AY55").Value = "NaN"
$ws.Range("W87").Value = 67
$ws.Range("W175").Value = 980

# --- Minor numeric corrections in column W (rows 167-194) ---
$ws.Range("W167").Value = 654
$ws.Range("W168").Value = 680
$ws.Range("W169").Value = 711
$ws.Range("W170").Value = 798
$ws.Range("W171").Value = 820
$ws.Range("W172").Value = 864
$ws.Range("W173").Value = 936
$ws.Range("W174").Value = 967
$ws.Range("W176").Value = 995
$ws.Range("W177").Value = 1029
$ws.Range("W178").Value = 1080
$ws.Range("W179").Value = 1211
$ws.Range("W180").Value = 1273
$ws.Range("W181").Value = 1303
$ws.Range("W182").Value = 1364
$ws.Range("W183").Value = 1392
$ws.Range("W184").Value = 1469
$ws.Range("W185").Value = 1527
$ws.Range("W186").Value = 1570
$ws.Range("W187").Value = 1615
$ws.Range("W189").Value = 1721
$ws.Range("W190").Value = 1780
$ws.Range("W192").Value = 1970
$ws.Range("W193").Value = 2027
$ws.Range("W194").Value = 2093

# --- Append new row 201 with the day's data ---
$ws.Range("A201").Value = 44095
$ws.Range("B201").Value = 770435
$ws.Range("C201").Value = 2732
$ws.Range("D201").Value = 104788
$ws.Range("E201").Value = 66643
$ws.Range("F201").Value = 254274
$ws.Range("G201").Value = 28058
$ws.Range("H201").Value = 5963
$ws.Range("I201").Value = 4785
$ws.Range("J201").Value = 7789
$ws.Range("K201").Value = 8340
$ws.Range("L201").Value = 17452
$ws.Range("M201").Value = 3958
$ws.Range("N201").Value = 22874
$ws.Range("O201").Value = 30838
$ws.Range("P201").Value = 7479
$ws.Range("Q201").Value = 9397
$ws.Range("R201").Value = 14572
$ws.Range("S201").Value = 13375
$ws.Range("T201").Value = 17556
$ws.Range("U201").Value = 14570
$ws.Range("V201").Value = 3652
$ws.Range("W201").Value = 2728
$ws.Range("X201").Value = 9419
$ws.Range("Y201").Value = 27835
$ws.Range("Z201").Value = 13619
$ws.Range("AA201").Value = 10969
$ws.Range("AB201").Value = 57483
$ws.Range("AC201").Value = 1929
$ws.Range("AD201").Value = 980
$ws.Range("AE201").Value = 693
$ws.Range("AF201").Value = 468
$ws.Range("AG201").Value = 638
$ws.Range("AH201").Value = 451
$ws.Range("AI201").Value = 624
$ws.Range("AJ201").Value = 2025
$ws.Range("AK201").Value = 5118
$ws.Range("AL201").Value = 37626
$ws.Range("AM201").Value = 9142
$ws.Range("AN201").Value = 2538
$ws.Range("AO201").Value = 44722
$ws.Range("AP201").Value = 1081
$ws.Range("AQ201").Value = 22444
$ws.Range("AR201").Value = 1519
$ws.Range("AS201").Value = 9948
$ws.Range("AT201").Value = 1644
$ws.Range("AU201").Value = 1602
$ws.Range("AV201").Value = 7514
$ws.Range("AW201").Value = 1976
$ws.Range("AX201").Value = 955
$ws.Range("AY201").Value = 2491
$ws.Range("AZ201").Value = 2662
$ws.Range("BA201").Value = 60863
$ws.Range("BB201").Value = 13697
$ws.Range("BC201").Value = 5639
$ws.Range("BD201").Value = 9477
$ws.Range("BE201").Value = 6319
$ws.Range("BF201").Value = 277
$ws.Range("BG201").Value = 1456
$ws.Range("BH201").Value = 2713
$ws.Range("BI201").Value = 741
$ws.Range("BJ201").Value = 2142
$ws.Range("BK201").Value = 9549
$ws.Range("BL201").Value = 9395
$ws.Range("BM201").Value = 10283
$ws.Range("BN201").Value = 14227
$ws.Range("BO201").Value = 1961
$ws.Range("BP201").Value = 894
$ws.Range("BQ201").Value = 12827
$ws.Range("BR201").Value = 10412
$ws.Range("BS201").Value = 12266
$ws.Range("BT201").Value = 2565
$ws.Range("BU201").Value = 2045
$ws.Range("BV201").Value = 5291
$ws.Range("BW201").Value = 4545
$ws.Range("BX201").Value = 1962
$ws.Range("BY201").Value = 5671
$ws.Range("BZ201").Value = 3392
$ws.Range("CA201").Value = 1969
$ws.Range("CB201").Value = 921
$ws.Range("CC201").Value = 2872
$ws.Range("CD201").Value = 2210
$ws.Range("CE201").Value = 1837
$ws.Range("CF201").Value = 1539
$ws.Range("CG201").Value = 6035
$ws.Range("CH201").Value = 2027
$ws.Range("CI201").Value = 1415
$ws.Range("CJ201").Value = 1734
$ws.Range("CK201").Value = 2053
$ws.Range("CL201").Value = 2066
$ws.Range("CM201").Value = 2457
$ws.Range("CN201").Value = 1640
$ws.Range("CO201").Value = 1205
$ws.Range("CP201").Value = 1197
$ws.Range("CQ201").Value = 923
$ws.Range("CR201").Value = 3343
$ws.Range("CS201").Value = 1422
$ws.Range("CT201").Value = 932
$ws.Range("CU201").Value = 1016
$ws.Range("CV201").Value = 1689
$ws.Range("CW201").Value = 1547
$ws.Range("CX201").Value = 751
$ws.Range("CY201").Value = 859
$ws.Range("CZ201").Value = 1267
$ws.Range("DA201").Value = 1566
$ws.Range("DB201").Value = 1468
$ws.Range("DC201").Value = 1512
$ws.Range("DD201").Value = 1185
$ws.Range("DE201").Value = 334
$ws.Range("DF201").Value = 365
$ws.Range("DG201").Value = 793
$ws.Range("DH201").Value = 752
$ws.Range("DI201").Value = 478
$ws.Range("DJ201").Value = 543
$ws.Range("DK201").Value = 380
$ws.Range("DL201").Value = 665
$ws.Range("DM201").Value = 749
$ws.Range("DN201").Value = 526
$ws.Range("DO201").Value = 490
$ws.Range("DP201").Value = 373
$ws.Range("DQ201").Value = 521
$ws.Range("DR201").Value = 134204
$ws.Range("DS201").Value = 326537
$ws.Range("DT201").Value = 17205
$ws.Range("DU201").Value = 141271
$ws.Range("DV201").Value = 87337
$ws.Range("DW201").Value = 42454
$ws.Range("DX201").Value = 12079

# --- Update active selection to reflect the new last cell ---
[void]$ws.Range("DX201").Select()
